$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header updates
# ---------------------------------------------------------------------------
# "VALOR MORA" total
$ws.Range("E11").Value = 16185960
# "Cant. Periodos" count (one new period added)
$ws.Range("F13").Value = 90

# ---------------------------------------------------------------------------
# 2) Insert a new row for the newest period (2508) right after the current
#    last data row (104), pushing the trailing "firma" block down by one row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(105).Insert()

# Row 105 inherits the "bottom of table" border styling that used to live on
# row 104; row 104 becomes a normal interior row. Do this by copying the
# formatting only (not the values) between rows.
$ws.Range("B104:J104").Copy()
$ws.Range("B105:J105").PasteSpecial(-4122)

$ws.Range("B103:J103").Copy()
$ws.Range("B104:J104").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Re-sort the "Periodo Mora" column: the table used to list periods from
#    newest (2507) to oldest (1803); it now lists them oldest (1803) to
#    newest (2508), the newest 2508 row having just been inserted above.
# ---------------------------------------------------------------------------
$periods = @("1803","1804","1805","1806","1807","1808","1809","1810","1811","1812", `
             "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912", `
             "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012", `
             "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112", `
             "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212", `
             "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312", `
             "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412", `
             "2501","2502","2503","2504","2505","2506","2507","2508")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = "1047369681"
    $ws.Cells.Item($row, 4).Value = "JULIAN ERNESTO DIAZ ARBOLEDA"
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = 179844
    $ws.Cells.Item($row, 7).Value = 4496100
}
